$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right count B11 5 -> 4, Wrong marking C11 -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total B12 65 -> 52, Max text E12 "65 / 140" -> "52 / 112"
$ws.Range("B12").Value = 52
$ws.Range("E12").Value = "52 / 112"
